$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value (matches the source inlineStr cells)
    # instead of Excel auto-coercing numeric-looking strings into numbers,
    # then restore the default "Normal" style so no stray number format sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Update Price (D) and Volume(1h) (E) columns ---
Set-TextValue $ws.Range("D2") "28.359.19"
$ws.Range("E2").Value = "  +4.82%  "
Set-TextValue $ws.Range("D3") "1.810.86"
$ws.Range("E3").Value = "  +3.93%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  +0.34%  "
Set-TextValue $ws.Range("D5") "317.85"
$ws.Range("E5").Value = "  +2.37%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  +0.37%  "
Set-TextValue $ws.Range("D7") "0.5722"
$ws.Range("E7").Value = "  +15.64%  "
Set-TextValue $ws.Range("D8") "0.3891"
$ws.Range("E8").Value = "  +10.46%  "
Set-TextValue $ws.Range("D9") "0.07596"
$ws.Range("E9").Value = "  +4.44%  "
Set-TextValue $ws.Range("D10") "42.91"
$ws.Range("E10").Value = "  +0.07%  "
Set-TextValue $ws.Range("D11") "1.139"
$ws.Range("E11").Value = "  +7.34%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("E13").Value = "  +5.46%  "
Set-TextValue $ws.Range("D14") "6.269"
$ws.Range("E14").Value = "  +6.16%  "
Set-TextValue $ws.Range("D15") "1.808.01"
$ws.Range("E15").Value = "  +4.44%  "
Set-TextValue $ws.Range("D16") "7.285"
$ws.Range("E16").Value = "  +6.44%  "
Set-TextValue $ws.Range("D17") "91.98"
$ws.Range("E17").Value = "  +5.33%  "
Set-TextValue $ws.Range("D18") "0.00001073"
$ws.Range("E18").Value = "  +3.46%  "
Set-TextValue $ws.Range("D19") "0.06477"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("E20").Value = "  +0.31%  "
Set-TextValue $ws.Range("D21") "17.30"
$ws.Range("E21").Value = "  +4.01%  "
Set-TextValue $ws.Range("D22") "6.004"
$ws.Range("E22").Value = "  +4.52%  "
Set-TextValue $ws.Range("D23") "28.367.71"
$ws.Range("E23").Value = "  +4.71%  "
Set-TextValue $ws.Range("D24") "11.32"
$ws.Range("E24").Value = "  +3.04%  "
Set-TextValue $ws.Range("D25") "2.137"
$ws.Range("E25").Value = "  +3.75%  "
Set-TextValue $ws.Range("D28") "20.78"
$ws.Range("E28").Value = "  +3.46%  "
Set-TextValue $ws.Range("D29") "2.019.96"
$ws.Range("E29").Value = "  +4.60%  "
Set-TextValue $ws.Range("D30") "124.25"
$ws.Range("E30").Value = "  +3.01%  "
Set-TextValue $ws.Range("D31") "1.166"
$ws.Range("E31").Value = "  +10.50%  "
Set-TextValue $ws.Range("D32") "0.1069"
$ws.Range("E32").Value = "  +13.91%  "
Set-TextValue $ws.Range("D33") "5.795"
$ws.Range("E33").Value = "  +7.04%  "
Set-TextValue $ws.Range("D34") "3.631"
$ws.Range("E34").Value = "  +1.39%  "
Set-TextValue $ws.Range("D35") "0.2211"
$ws.Range("E35").Value = "  +10.80%  "
Set-TextValue $ws.Range("D36") "8.976"
$ws.Range("E36").Value = "  +20.47%  "
$ws.Range("E37").Value = "  +5.79%  "
Set-TextValue $ws.Range("D38") "11.68"
$ws.Range("E38").Value = "  +5.79%  "
Set-TextValue $ws.Range("D39") "0.06132"
$ws.Range("E39").Value = "  +3.16%  "
Set-TextValue $ws.Range("D40") "0.6399"
$ws.Range("E40").Value = "  +6.01%  "
Set-TextValue $ws.Range("D41") "5.036"
$ws.Range("E41").Value = "  +5.66%  "
Set-TextValue $ws.Range("D42") "1.164"
$ws.Range("E42").Value = "  +4.36%  "
$ws.Range("E43").Value = "  +0.29%  "
Set-TextValue $ws.Range("D44") "1.380"
$ws.Range("E44").Value = "  -3.59%  "
Set-TextValue $ws.Range("D45") "13.49"
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("E46").Value = "  +6.44%  "
Set-TextValue $ws.Range("D47") "3.696"
$ws.Range("E47").Value = "  +3.37%  "
Set-TextValue $ws.Range("D48") "122.89"
$ws.Range("E48").Value = "  +2.63%  "
Set-TextValue $ws.Range("D49") "1.957"
$ws.Range("E49").Value = "  +5.77%  "
$ws.Range("E50").Value = "  +4.33%  "
Set-TextValue $ws.Range("D51") "0.06877"
$ws.Range("E51").Value = "  +2.89%  "

# --- Rows 26 and 27 swap content (LidoDAOToken <-> Monero) with updated Price/Volume ---
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "158.50"
$ws.Range("E26").Value = "  +3.10%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D27") "2.461"
$ws.Range("E27").Value = "  +17.56%  "
